$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text so values like "1.009" or
# "0.000008590" are not reinterpreted as numbers (losing formatting/precision).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.543.06'
$ws.Range("E2").Value = '  -2.61%  '

# Row 3
$ws.Range("D3").Value = '1.809.37'
$ws.Range("E3").Value = '  -2.34%  '

# Row 4
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.76%  '

# Row 5
$ws.Range("D5").Value = '1.008'
$ws.Range("E5").Value = '  +0.65%  '

# Row 6
$ws.Range("D6").Value = '308.86'
$ws.Range("E6").Value = '  -1.66%  '

# Row 7
$ws.Range("D7").Value = '0.4545'
$ws.Range("E7").Value = '  -1.30%  '

# Row 8
$ws.Range("D8").Value = '0.3659'
$ws.Range("E8").Value = '  -1.29%  '

# Row 9
$ws.Range("D9").Value = '0.07122'
$ws.Range("E9").Value = '  -2.43%  '

# Row 10
$ws.Range("D10").Value = '0.8767'
$ws.Range("E10").Value = '  -1.21%  '

# Row 11
$ws.Range("D11").Value = '0.07758'
$ws.Range("E11").Value = '  -0.61%  '

# Row 12
$ws.Range("D12").Value = '19.41'
$ws.Range("E12").Value = '  -3.38%  '

# Row 13
$ws.Range("D13").Value = '1.801.93'
$ws.Range("E13").Value = '  -4.73%  '

# Row 14
$ws.Range("D14").Value = '5.270'
$ws.Range("E14").Value = '  -2.19%  '

# Row 15
$ws.Range("D15").Value = '6.342'
$ws.Range("E15").Value = '  -2.86%  '

# Row 16
$ws.Range("D16").Value = '86.29'
$ws.Range("E16").Value = '  -5.62%  '

# Row 17
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  +0.78%  '

# Row 18
$ws.Range("D18").Value = '0.000008590'
$ws.Range("E18").Value = '  -3.87%  '

# Row 19
$ws.Range("D19").Value = '1.008'
$ws.Range("E19").Value = '  +0.63%  '

# Row 20
$ws.Range("D20").Value = '26.597.89'
$ws.Range("E20").Value = '  -2.48%  '

# Row 21
$ws.Range("D21").Value = '14.26'
$ws.Range("E21").Value = '  -3.41%  '

# Row 22
$ws.Range("D22").Value = '4.971'
$ws.Range("E22").Value = '  -2.71%  '

# Row 23
$ws.Range("D23").Value = '10.39'
$ws.Range("E23").Value = '  -1.38%  '

# Row 24
$ws.Range("D24").Value = '1.972'
$ws.Range("E24").Value = '  +2.23%  '

# Row 25
$ws.Range("D25").Value = '150.86'
$ws.Range("E25").Value = '  -0.64%  '

# Row 26
$ws.Range("D26").Value = '17.92'
$ws.Range("E26").Value = '  -2.82%  '

# Row 27
$ws.Range("D27").Value = '2.004'
$ws.Range("E27").Value = '  -2.61%  '

# Row 28
$ws.Range("D28").Value = '113.01'
$ws.Range("E28").Value = '  -2.45%  '

# Row 29
$ws.Range("D29").Value = '4.858'
$ws.Range("E29").Value = '  -4.02%  '

# Row 30
$ws.Range("D30").Value = '0.08669'
$ws.Range("E30").Value = '  -1.78%  '

# Row 31
$ws.Range("D31").Value = '3.062'
$ws.Range("E31").Value = '  -1.15%  '

# Row 32
$ws.Range("D32").Value = '4.458'
$ws.Range("E32").Value = '  -0.96%  '

# Row 33
$ws.Range("D33").Value = '0.7291'
$ws.Range("E33").Value = '  -5.41%  '

# Row 34
$ws.Range("D34").Value = '1.115'
$ws.Range("E34").Value = '  -4.53%  '

# Row 35
$ws.Range("D35").Value = '1.006'
$ws.Range("E35").Value = '  +0.48%  '

# Row 36
$ws.Range("D36").Value = '2.555'
$ws.Range("E36").Value = '  -7.14%  '

# Row 37
$ws.Range("D37").Value = '1.082'
$ws.Range("E37").Value = '  -0.15%  '

# Row 38
$ws.Range("D38").Value = '0.01929'
$ws.Range("E38").Value = '  -1.06%  '

# Row 39
$ws.Range("D39").Value = '0.05088'
$ws.Range("E39").Value = '  -3.25%  '

# Row 40
$ws.Range("D40").Value = '2.875'
$ws.Range("E40").Value = '  -2.43%  '

# Row 41
$ws.Range("D41").Value = '6.959'
$ws.Range("E41").Value = '  -1.54%  '

# Row 42
$ws.Range("D42").Value = '0.4977'
$ws.Range("E42").Value = '  -2.82%  '

# Row 43
$ws.Range("D43").Value = '0.1568'
$ws.Range("E43").Value = '  -3.99%  '

# Row 44
$ws.Range("D44").Value = '8.133'
$ws.Range("E44").Value = '  -3.05%  '

# Row 45
$ws.Range("D45").Value = '1.009'
$ws.Range("E45").Value = '  +0.76%  '

# Row 46
$ws.Range("D46").Value = '0.4622'
$ws.Range("E46").Value = '  -3.66%  '

# Row 47
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '101.75'
$ws.Range("E47").Value = '  -0.51%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.912'
$ws.Range("E48").Value = '  -3.82%  '

# Row 49
$ws.Range("D49").Value = '1.587'
$ws.Range("E49").Value = '  -3.36%  '

# Row 50
$ws.Range("D50").Value = '0.06001'
$ws.Range("E50").Value = '  -3.41%  '

# Row 51
$ws.Range("D51").Value = '63.78'
$ws.Range("E51").Value = '  -2.90%  '
